# Update code list dictionaries (roles.xlsx)
# - Rename several role IDs/names (checked_by->checker, prepared_by->preparer,
#   calculated_by->calculator, tested_by->tester) and refresh their descriptions.
# - Add a new "collector" role to the Definitions table.
# - Update SourceElement XPaths in AssociatedElements to the new DIGGS paths
#   and add several new AssociatedElements rows covering the new/renamed roles.

$wb = $excel.ActiveWorkbook

$defs = $wb.Worksheets.Item("Definitions")
$assoc = $wb.Worksheets.Item("AssociatedElements")

$defsTbl = $defs.ListObjects.Item("Definitions")
$assocTbl = $assoc.ListObjects.Item("AssociatedElements")

# ---------------------------------------------------------------------------
# Definitions sheet updates
# ---------------------------------------------------------------------------

# Row 8 (Project Manager) - value re-entered (no textual change)
$defs.Range("B8").Value = "project_manager"

# Row 9: Checked By -> Checker
$defs.Range("B9").Value = "checker"
$defs.Range("C9").Value = "Checker"
$defs.Range("D9").Value = "Person checking the associated data, or for accuracy or standards compliance"

# Row 10: Prepared By -> Preparer
$defs.Range("B10").Value = "preparer"
$defs.Range("C10").Value = "Preparer"
$defs.Range("D10").Value = "Person compiling or inputting the data for  or ppreparing specimen for testing"

# Row 11: Calculated By -> Calculator
$defs.Range("B11").Value = "calculator"
$defs.Range("C11").Value = "Calculator"
$defs.Range("D11").Value = "Person performing calculations to obtain derived results"

# Row 12: Tested By -> Tester (description unchanged)
$defs.Range("C12").Value = "Tester"

# Row 15: Helper description updated
$defs.Range("D15").Value = "Person assisting with construction (eg. borehole, actvity, sampling, observation or measurement"

# New row 16: Collector
$defsNewRow = $defsTbl.ListRows.Add()
$defs.Range("A16").Formula = '=IF(ISNA(VLOOKUP(B16,AssociatedElements!B$2:B2850,1,FALSE)),"Not used","")'
$defs.Range("B16").Value = "collector"
$defs.Range("C16").Value = "Collector"
$defs.Range("D16").Value = "Person that collects or creates a sample"
$defs.Range("E16").Value = "string"
$defs.Range("G16").Value = "DIGGS"

# ---------------------------------------------------------------------------
# AssociatedElements sheet updates
# ---------------------------------------------------------------------------

# Existing rows: refresh SourceElement Xpaths and renamed IDs
$assoc.Range("C2").Value = "//diggs:samplingFeature//rolePerformed"
$assoc.Range("C3").Value = "//diggs:samplingFeature//rolePerformed"
$assoc.Range("C4").Value = "//diggs:samplingFeature//rolePerformed"

$assoc.Range("C6").Value = "/diggs:measurement//diggs:rolePerformed"

$assoc.Range("C7").Value = "//diggs:project//diggs:rolePerformed"
$assoc.Range("C8").Value = "//diggs:project//diggs:rolePerformed"

$assoc.Range("B9").Value = "checker"
$assoc.Range("C9").Value = "//diggs:rolePerformed"

$assoc.Range("B10").Value = "preparer"

$assoc.Range("B11").Value = "calculator"

$assoc.Range("B12").Value = "tester"

$assoc.Range("C13").Value = "//diggs:project//diggs:rolePerformed"
$assoc.Range("C14").Value = "//diggs:project//diggs:rolePerformed"

$assoc.Range("C15").Value = "//diggs:samplingFeature//rolePerformed"

# New rows 16-22
for ($i = 0; $i -lt 7; $i++) {
    [void]$assocTbl.ListRows.Add()
}

$assoc.Range("B16").Value = "operator"
$assoc.Range("C16").Value = "//diggs:samplingActivity//rolePerformed"

$assoc.Range("B17").Value = "logger"
$assoc.Range("C17").Value = "//diggs:rolePerformed"

$assoc.Range("B18").Value = "drilling_contractor"
$assoc.Range("C18").Value = "//diggs:samplingActivity//rolePerformed"

$assoc.Range("B19").Value = "helper"
$assoc.Range("C19").Value = "//diggs:rolePerformed"

$assoc.Range("B20").Value = "collector"
$assoc.Range("C20").Value = "//diggs:samplingActivity//rolePerformed"

$assoc.Range("B21").Value = "collector"
$assoc.Range("C21").Value = "//diggs:sample//rolePerformed"

$assoc.Range("B22").Value = "laboratory_name"
$assoc.Range("C22").Value = "//diggs:project//diggs:rolePerformed"
